$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ignitions")

# Fill in C column (Feet) for HFTD2 rows (8,9,10)
$ws.Range("C8").Value = 85018055
$ws.Range("C9").Value = 11977176
$ws.Range("C10").Value = 3697017

# Fill in C column (Feet) for HFTD3 rows (13,14,15)
$ws.Range("C13").Value = 30827591
$ws.Range("C14").Value = 5324275
$ws.Range("C15").Value = 3416580

# C18:C20 sum formulas (HFTD = HFTD2 + HFTD3)
$ws.Range("C18").Formula = "=C8+C13"
$ws.Range("C19").Formula = "=C9+C14"
$ws.Range("C20").Formula = "=C10+C15"

# D column formulas (Miles = Feet/5280), one shared formula per block
$ws.Range("D8:D10").Formula = "=C8/5280"
$ws.Range("D13:D15").Formula = "=C13/5280"
$ws.Range("D18:D20").Formula = "=C18/5280"

# E column formulas (Ign/Mile = Ignitions/Miles), one shared formula per block
$ws.Range("E8:E10").Formula = "=B8/D8"
$ws.Range("E13:E15").Formula = "=B13/D13"
$ws.Range("E18:E20").Formula = "=B18/D18"

# Update the selection on the sheet
$ws.Range("A1:E20").Select()
